# Generate Report for Handback
# Updates the localization-status workbook to reflect a handback
# transform failure for the 3dad73a3-... file, across the Overview,
# zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusText = "Handback transform failed"

$zhError = "Handback file name: 0si3h5ag.fvc is different with handoff file name: 3dad73a3-6d56-4070-8073-1543cc1217bf.f4f54412e95a4ad77398d9e3f84798bb74455cd9.zh-cn."
$deError  = "Handback file name: 0si3h5ag.fvc is different with handoff file name: 3dad73a3-6d56-4070-8073-1543cc1217bf.f4f54412e95a4ad77398d9e3f84798bb74455cd9.de-de."

# --- Overview sheet: update status for the 3dad73a3... row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# Excel's ColumnWidth property is offset from the raw OOXML column
# width (in characters) by the default column padding (5px at the
# workbook's default Calibri 11 font == 0.8333... characters). Setting
# ColumnWidth = 39.1666... yields a stored width of exactly 40.
$targetWidth = 39.166666666666664

# --- zh-cn sheet: update status, error detail, and column width ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = $statusText
$wsZh.Range("P3").Value = $zhError
$wsZh.Columns.Item(16).ColumnWidth = $targetWidth

# --- de-de sheet: update status, error detail, and column width ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = $statusText
$wsDe.Range("P3").Value = $deError
$wsDe.Columns.Item(16).ColumnWidth = $targetWidth
